$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Key" item category (columns Z/AA), headed by "Key", row 3 holds "Iron Key" / 7001
$ws.Range("Z1").Value = "Key"
$ws.Range("Z3").Value = "Iron Key"
$ws.Range("AA3").Value = 7001

# New "Fireball" spell entry added under the existing Spell column (Q/R), row 6
$ws.Range("Q6").Value = "Fireball"
$ws.Range("R6").Value = 4004

# Match the author's final selection/viewport on the new data
$ws.Range("R6").Select()
